$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Issue #268: ISBLANK() on a formula-produced empty string should still
# be treated as blank. Add a regression test row below the existing data.
$ws.Range("A16").Formula = '=IF(ISBLANK(B16), "Düsseldorf", B16)'

$ws.Range("A16").Select()

